$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several Price (column D) values look like plain numbers (e.g. "1.001"),
# but in the source workbook every Price/Volume cell is stored as literal text
# (e.g. "25.798.77" with two dots is clearly not numeric). To keep that text
# typing for the values that *would* otherwise be auto-parsed as a number by
# Excel, those assignments are prefixed with a leading apostrophe, exactly as
# a user forcing text entry would type them in the UI.

$ws.Range('D2').Value = '25.859.29'
$ws.Range('E2').Value = '  -0.25%  '

$ws.Range('D3').Value = '1.630.91'
$ws.Range('E3').Value = '  -0.50%  '

$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '''215.63'
$ws.Range('E5').Value = '  +0.42%  '

$ws.Range('D6').Value = '''0.5063'
$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').Value = '''1.001'
$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').Value = '''0.2570'
$ws.Range('E8').Value = '  +0.69%  '

$ws.Range('D9').Value = '''0.06336'
$ws.Range('E9').Value = '  -0.44%  '

$ws.Range('D10').Value = '''19.49'
$ws.Range('E10').Value = '  +0.11%  '

$ws.Range('D11').Value = '''0.07755'
$ws.Range('E11').Value = '  +0.34%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.650.35'
$ws.Range('E12').Value = '  -0.05%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.240'
$ws.Range('E13').Value = '  -0.81%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.853.10'
$ws.Range('E14').Value = '  -0.77%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.5484'
$ws.Range('E15').Value = '  +0.86%  '

$ws.Range('D16').Value = '0.0₅7656'
$ws.Range('E16').Value = '  -1.83%  '

$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '''63.67'
$ws.Range('E17').Value = '  -0.78%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '25.857.76'
$ws.Range('E18').Value = '  -0.35%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '''1.002'
$ws.Range('E19').Value = '  +0.01%  '

$ws.Range('D20').Value = '''194.94'
$ws.Range('E20').Value = '  -0.89%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '''4.423'
$ws.Range('E21').Value = '  -0.75%  '

$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').Value = '''9.889'
$ws.Range('E22').Value = '  -0.39%  '

$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = '''6.059'
$ws.Range('E23').Value = '  +0.70%  '

$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').Value = '''1.002'
$ws.Range('E24').Value = '  -0.30%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''1.923'
$ws.Range('E25').Value = '  +1.92%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''141.97'
$ws.Range('E26').Value = '  +0.73%  '

$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '''0.1233'
$ws.Range('E27').Value = '  +3.95%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''6.788'
$ws.Range('E28').Value = '  -1.13%  '

$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''15.57'
$ws.Range('E29').Value = '  -0.69%  '

$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.238'
$ws.Range('E30').Value = '  +0.21%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.04872'
$ws.Range('E31').Value = '  -1.26%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '''3.237'
$ws.Range('E32').Value = '  -0.45%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''3.181'
$ws.Range('E33').Value = '  +0.19%  '

$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '''1.537'
$ws.Range('E34').Value = '  +0.02%  '

$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.371'
$ws.Range('E35').Value = '  +0.10%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''0.8976'
$ws.Range('E36').Value = '  +0.48%  '

$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '''2.540'
$ws.Range('E37').Value = '  -1.47%  '

$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '''0.5499'
$ws.Range('E38').Value = '  +1.46%  '

$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.117.30'
$ws.Range('E39').Value = '  -1.34%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '''0.01552'
$ws.Range('E40').Value = '  -0.08%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '''1.001'
$ws.Range('E41').Value = '  -0.13%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '''5.561'
$ws.Range('E42').Value = '  -0.29%  '

$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '''0.7968'
$ws.Range('E43').Value = '  -2.25%  '

$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '''97.18'
$ws.Range('E44').Value = '  -2.11%  '

$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₈118'
$ws.Range('E45').Value = '  -8.37%  '

$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D46').Value = '1.761.21'
$ws.Range('E46').Value = '  -0.85%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '''0.4438'
$ws.Range('E47').Value = '  -2.15%  '

$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = '''1.001'
$ws.Range('E48').Value = '  -0.21%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '''54.64'
$ws.Range('E49').Value = '  -0.19%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.05138'
$ws.Range('E50').Value = '  +1.28%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''7.543'
$ws.Range('E51').Value = '  +2.84%  '
